# Update "想去人数" (F column) counts across sheets to reflect the latest
# generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 503
$ws.Range("F7").Value = 1126
$ws.Range("F8").Value = 725
$ws.Range("F10").Value = 1396
$ws.Range("F11").Value = 276
$ws.Range("F14").Value = 63
$ws.Range("F16").Value = 45
$ws.Range("F18").Value = 14
$ws.Range("F20").Value = 292
$ws.Range("F21").Value = 545
$ws.Range("F22").Value = 559
$ws.Range("F23").Value = 751

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 992
$ws.Range("F5").Value = 257

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 213

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 503
$ws.Range("F3").Value = 213
$ws.Range("F8").Value = 992
$ws.Range("F10").Value = 1126
$ws.Range("F11").Value = 725
$ws.Range("F13").Value = 1396
$ws.Range("F14").Value = 276
$ws.Range("F17").Value = 63
$ws.Range("F19").Value = 45
$ws.Range("F21").Value = 14
$ws.Range("F23").Value = 257
$ws.Range("F25").Value = 292
$ws.Range("F29").Value = 545
$ws.Range("F30").Value = 559
$ws.Range("F31").Value = 751
